$d = $word.ActiveDocument

# --- First paragraph: "**ID__AFFARS_5337_topic_7__ID** " identifier line ---
$p1 = $d.Paragraphs(1)

# Add a 5-twip paragraph border on all four sides.
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Replace the placeholder id text and drop the trailing space-only run,
# merging everything into a single run with the new id text.
$r1 = $p1.Range
$r1.End = $r1.End - 1
$r1.Text = "**ID__AFFARS_SUBPART_5337_2__ID**"

Write-Output "done"
